$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Vtn-Itgb3 LR-pairs data ("Natmi following Dr Hou advice"):
# sending/target clusters now include the new "M2" cluster (4x4=16 combinations,
# rows 2-17) and the ligand/receptor summary stats were recomputed accordingly.
$rows = @(
    @("ECs", "Vtn", "Itgb3", "ECs", 3, 1, 12.75206033333333, 38.256181, 0.1573122343381959, 0.157312234338196, 3, 1, 3.778439, 11.335317, 0.4252971528324392, 0.4252971528324392, 48.18288209381966, 433.645938844377, 0.06690444536974421, 0.06690444536974423),
    @("ECs", "Vtn", "Itgb3", "FAPs", 3, 1, 12.75206033333333, 38.256181, 0.1573122343381959, 0.157312234338196, 3, 1, 4.333403333333333, 13.00021, 0.4877633593505858, 0.4877633593505858, 55.25982075533444, 497.3383867980099, 0.07673114388774503, 0.07673114388774505),
    @("ECs", "Vtn", "Itgb3", "sCs", 3, 1, 12.75206033333333, 38.256181, 0.1573122343381959, 0.157312234338196, 3, 1, 0.2909853333333334, 0.8729560000000001, 0.03275300561492853, 0.03275300561492853, 3.710662526781778, 33.395962741036, 0.005152448494575885, 0.005152448494575886),
    @("ECs", "Vtn", "Itgb3", "M2", 3, 1, 12.75206033333333, 38.256181, 0.1573122343381959, 0.157312234338196, 3, 1, 0.4814053333333333, 1.444216, 0.0541864822020464, 0.0541864822020464, 6.13890985545511, 55.25018869909599, 0.008524196586130807, 0.00852419658613081),
    @("FAPs", "Vtn", "Itgb3", "ECs", 3, 1, 20.35396833333334, 61.06190500000001, 0.2510910513649196, 0.2510910513649196, 3, 1, 3.778439, 11.335317, 0.4252971528324392, 0.4252971528324392, 76.90622775543169, 692.1560497988851, 0.1067883092472041, 0.1067883092472041),
    @("FAPs", "Vtn", "Itgb3", "FAPs", 3, 1, 20.35396833333334, 61.06190500000001, 0.2510910513649196, 0.2510910513649196, 3, 1, 4.333403333333333, 13.00021, 0.4877633593505858, 0.4877633593505858, 88.2019542222278, 793.81758800005, 0.1224730147166237, 0.1224730147166237),
    @("FAPs", "Vtn", "Itgb3", "sCs", 3, 1, 20.35396833333334, 61.06190500000001, 0.2510910513649196, 0.2510910513649196, 3, 1, 0.2909853333333334, 0.8729560000000001, 0.03275300561492853, 0.03275300561492853, 5.922706260131113, 53.30435634118001, 0.00822398661521352, 0.00822398661521352),
    @("FAPs", "Vtn", "Itgb3", "M2", 3, 1, 20.35396833333334, 61.06190500000001, 0.2510910513649196, 0.2510910513649196, 3, 1, 0.4814053333333333, 1.444216, 0.0541864822020464, 0.05418648220204641, 9.798508910164445, 88.18658019148, 0.01360574078587833, 0.01360574078587834),
    @("sCs", "Vtn", "Itgb3", "ECs", 1, 0.3333333333333333, 0.004706, 0.014118, 0.000058054255974652841049230695, 0.000058054255974652847825494273, 3, 1, 3.778439, 11.335317, 0.4252971528324392, 0.4252971528324392, 0.017781333934, 0.160032005406, 0.000024690309775825480348388982, 0.000024690309775825480348388982),
    @("sCs", "Vtn", "Itgb3", "FAPs", 1, 0.3333333333333333, 0.004706, 0.014118, 0.000058054255974652841049230695, 0.000058054255974652847825494273, 3, 1, 4.333403333333333, 13.00021, 0.4877633593505858, 0.4877633593505858, 0.02039299608666667, 0.18353696478, 0.000028316738918795489522066758, 0.000028316738918795489522066758),
    @("sCs", "Vtn", "Itgb3", "sCs", 1, 0.3333333333333333, 0.004706, 0.014118, 0.000058054255974652841049230695, 0.000058054255974652847825494273, 3, 1, 0.2909853333333334, 0.8729560000000001, 0.03275300561492853, 0.03275300561492853, 0.001369376978666667, 0.012324392808, 0.000001901451371908303095427622, 0.000001901451371908303095427622),
    @("sCs", "Vtn", "Itgb3", "M2", 1, 0.3333333333333333, 0.004706, 0.014118, 0.000058054255974652841049230695, 0.000058054255974652847825494273, 3, 1, 0.4814053333333333, 1.444216, 0.0541864822020464, 0.05418648220204641, 0.002265493498666666, 0.020389441488, 0.000003145755908123572106753832, 0.000003145755908123572953786779),
    @("M2", "Vtn", "Itgb3", "ECs", 3, 1, 47.951367, 143.854101, 0.5915386600409097, 0.5915386600409098, 3, 1, 3.778439, 11.335317, 0.4252971528324392, 0.4252971528324392, 181.181315176113, 1630.631836585017, 0.2515797079057151, 0.2515797079057152),
    @("M2", "Vtn", "Itgb3", "FAPs", 3, 1, 47.951367, 143.854101, 0.5915386600409097, 0.5915386600409098, 3, 1, 4.333403333333333, 13.00021, 0.4877633593505858, 0.4877633593505858, 207.79261359569, 1870.13352236121, 0.2885308840072983, 0.2885308840072983),
    @("M2", "Vtn", "Itgb3", "sCs", 3, 1, 47.951367, 143.854101, 0.5915386600409097, 0.5915386600409098, 3, 1, 0.2909853333333334, 0.8729560000000001, 0.03275300561492853, 0.03275300561492853, 13.953144510284, 125.578300592556, 0.01937466905376722, 0.01937466905376722),
    @("M2", "Vtn", "Itgb3", "M2", 3, 1, 47.951367, 143.854101, 0.5915386600409097, 0.5915386600409098, 3, 1, 0.4814053333333333, 1.444216, 0.0541864822020464, 0.05418648220204641, 23.084043814424, 207.756394329816, 0.03205339907412913, 0.03205339907412914)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    $excelRow = $i + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $row[$c]
    }
}
